Set-StrictMode -Version Latest
$ErrorActionPreference = "Stop"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Row 10 is the "R30" rule row; column C is the "From" hour. Update
# the "From" value for R30 from 18 to 1.
$ws.Range("C10").Value = 1
